$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion message text (cell A1) ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$oldLine1 = "✅ 1000 Bs = 3.24 = 12317.15 pesos"
$newLine1 = "✅ 1000 Bs = 3.18 = 12116.87 pesos"
$oldLine2 = "✅ 12317.15 pesos = 3.22 = 959.65 Bs"
$newLine2 = "✅ 12116.87 pesos = 3.17 = 957.04 Bs"

$text = $ws1.Range("A1").Value2
$text = $text.Replace($oldLine1, $newLine1)
$text = $text.Replace($oldLine2, $newLine2)
$ws1.Range("A1").Value2 = $text

# --- Sheet "tasas": update rate values ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 313.987
$ws2.Range("O10").Value = 3804.54
$ws2.Range("N12").Value = 3818.5
$ws2.Range("O12").Value = 301.6
